$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.473.38"
$ws.Range("E2").Value = "  +3.60%  "
$ws.Range("D3").Value = "3.460.24"
$ws.Range("E3").Value = "  +3.37%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "'407.98"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").Value = "'130.64"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +16.79%  "
$ws.Range("D7").Value = "3.453.67"
$ws.Range("E7").Value = "  +3.34%  "
$ws.Range("D8").Value = "'0.596"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.07%  "
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "'0.687"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +8.55%  "
$ws.Range("D11").Value = "'0.127"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +28.46%  "
$ws.Range("D12").Value = "'42.56"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +6.34%  "
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "4.019.24"
$ws.Range("E14").Value = "  +4.16%  "
$ws.Range("D15").Value = "'8.68"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "'20.00"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.39%  "
$ws.Range("D17").Value = "3.463.76"
$ws.Range("E17").Value = "  +3.59%  "
$ws.Range("D18").Value = "62.539.58"
$ws.Range("E18").Value = "  +4.42%  "
$ws.Range("D19").Value = "'1.04"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").Value = "'10.82"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("E21").Value = "  +22.94%  "
$ws.Range("D22").Value = "'3.36"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "'82.46"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +9.02%  "
$ws.Range("D24").Value = "'13.08"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "'308.04"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("E26").Value = "  -3.72%  "
$ws.Range("D27").Value = "'30.21"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.48%  "
$ws.Range("D28").Value = "'8.25"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.78%  "
$ws.Range("D29").Value = "'7.74"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.55%  "
$ws.Range("D30").Value = "'0.179"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.25%  "
$ws.Range("D31").Value = "'4.38"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("E32").Value = "  +4.21%  "
$ws.Range("D33").Value = "'2.67"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").Value = "'11.87"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.68%  "
$ws.Range("D35").Value = "'43.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +8.29%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "'0.0491"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("D38").Value = "'52.54"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("D39").Value = "'3.56"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.84%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("D41").Value = "'2.98"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -7.23%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'137.97"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.126"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.14%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'1.98"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.25%  "
$ws.Range("D45").Value = "'17.47"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.46%  "
$ws.Range("D46").Value = "'3.95"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").Value = "'0.284"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").Value = "'22.26"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("D50").Value = "2.202.21"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "3.810.23"
$ws.Range("E51").Value = "  +4.14%  "
